# Scheduled-runner refresh of cached market/profit figures across the
# per-job "Sheets" (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Each block
# below updates the currentAveragePrice*/LevePrice*/LeveProfit* columns
# (H:N) for a handful of rows per sheet to the latest pulled values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1575.2222
$ws.Range("I12").Value = 3365
$ws.Range("J12").Value = 143.4
$ws.Range("K12").Value = 3365
$ws.Range("L12").Value = 143.4
$ws.Range("M12").Value = -3195
$ws.Range("N12").Value = -483.4

$ws.Range("H70").Value = 73278
$ws.Range("I70").Value = 1625
$ws.Range("J70").Value = 101939.2
$ws.Range("K70").Value = 4875
$ws.Range("L70").Value = 305817.6
$ws.Range("M70").Value = -4605
$ws.Range("N70").Value = -306357.6

$ws.Range("H73").Value = 73278
$ws.Range("I73").Value = 1625
$ws.Range("J73").Value = 101939.2
$ws.Range("K73").Value = 4875
$ws.Range("L73").Value = 305817.6
$ws.Range("M73").Value = -3939
$ws.Range("N73").Value = -307689.6

$ws.Range("H92").Value = 323.53845
$ws.Range("I92").Value = 326.75
$ws.Range("K92").Value = 326.75
$ws.Range("M92").Value = 921.25

$ws.Range("H132").Value = 4101.35
$ws.Range("I132").Value = 4078.1538
$ws.Range("K132").Value = 12234.4614
$ws.Range("M132").Value = -9704.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 6745.25
$ws.Range("I63").Value = 2722
$ws.Range("K63").Value = 2722
$ws.Range("M63").Value = -2036

$ws.Range("H66").Value = 6745.25
$ws.Range("I66").Value = 2722
$ws.Range("K66").Value = 13610
$ws.Range("M66").Value = -10178

$ws.Range("H117").Value = 50248
$ws.Range("J117").Value = 50248
$ws.Range("L117").Value = 50248
$ws.Range("N117").Value = -59426

$ws.Range("H122").Value = 50006508
$ws.Range("I122").Value = 111116530
$ws.Range("J122").Value = 7403.636
$ws.Range("K122").Value = 333349590
$ws.Range("L122").Value = 22210.908
$ws.Range("M122").Value = -333347140
$ws.Range("N122").Value = -27110.908

$ws.Range("H132").Value = 238878.83
$ws.Range("I132").Value = 260379.3
$ws.Range("K132").Value = 781137.8999999999
$ws.Range("M132").Value = -778607.8999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 96428.664
$ws.Range("I134").Value = 4782.75
$ws.Range("K134").Value = 14348.25
$ws.Range("M134").Value = -11813.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 17050.691
$ws.Range("J41").Value = 20450
$ws.Range("L41").Value = 20450
$ws.Range("N41").Value = -21306

$ws.Range("H58").Value = 9975.210999999999
$ws.Range("I58").Value = 3635.889
$ws.Range("J58").Value = 15680.6
$ws.Range("K58").Value = 3635.889
$ws.Range("L58").Value = 15680.6
$ws.Range("M58").Value = -3432.889
$ws.Range("N58").Value = -16086.6

$ws.Range("H92").Value = 73999.5
$ws.Range("J92").Value = 73999.5
$ws.Range("L92").Value = 73999.5
$ws.Range("N92").Value = -78991.5

$ws.Range("I99").Value = 3244.5715
$ws.Range("J99").Value = 8166.5
$ws.Range("K99").Value = 3244.5715
$ws.Range("L99").Value = 8166.5
$ws.Range("M99").Value = -1746.5715
$ws.Range("N99").Value = -11162.5

$ws.Range("I126").Value = 3244.5715
$ws.Range("J126").Value = 8166.5
$ws.Range("K126").Value = 9733.7145
$ws.Range("L126").Value = 24499.5
$ws.Range("M126").Value = -7263.7145
$ws.Range("N126").Value = -29439.5

$ws.Range("H132").Value = 3398.1177
$ws.Range("I132").Value = 2126.2856
$ws.Range("K132").Value = 6378.8568
$ws.Range("M132").Value = -3848.8568

$ws.Range("H134").Value = 618233.5
$ws.Range("I134").Value = 41329.668
$ws.Range("J134").Value = 1267250.4
$ws.Range("K134").Value = 123989.004
$ws.Range("L134").Value = 3801751.2
$ws.Range("M134").Value = -121454.004
$ws.Range("N134").Value = -3806821.2

$ws.Range("H136").Value = 9975.210999999999
$ws.Range("I136").Value = 3635.889
$ws.Range("J136").Value = 15680.6
$ws.Range("K136").Value = 10907.667
$ws.Range("L136").Value = 47041.8
$ws.Range("M136").Value = -8357.667000000001
$ws.Range("N136").Value = -52141.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 4174.8
$ws.Range("I14").Value = 4174.8
$ws.Range("K14").Value = 12524.4
$ws.Range("M14").Value = -12351.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 411733.94
$ws.Range("I122").Value = 554468.4399999999
$ws.Range("J122").Value = 3921
$ws.Range("K122").Value = 1663405.32
$ws.Range("L122").Value = 11763
$ws.Range("M122").Value = -1660955.32
$ws.Range("N122").Value = -16663

$ws.Range("H132").Value = 48413.38
$ws.Range("I132").Value = 16623.262
$ws.Range("K132").Value = 49869.78599999999
$ws.Range("M132").Value = -47339.78599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 426772.4
$ws.Range("I7").Value = 719482.7
$ws.Range("J7").Value = 16978
$ws.Range("K7").Value = 719482.7
$ws.Range("L7").Value = 16978
$ws.Range("M7").Value = -719370.7
$ws.Range("N7").Value = -17202

$ws.Range("H46").Value = 2592.5278
$ws.Range("J46").Value = 2837.3125
$ws.Range("L46").Value = 2837.3125
$ws.Range("N46").Value = -3213.3125

$ws.Range("H93").Value = 1368.9412
$ws.Range("I93").Value = 1368.9412
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1368.9412
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -120.9412
$ws.Range("N93").ClearContents()

$ws.Range("H100").Value = 67405.53
$ws.Range("I100").Value = 82645.69500000001
$ws.Range("K100").Value = 82645.69500000001
$ws.Range("M100").Value = -82104.69500000001

$ws.Range("H122").Value = 528176
$ws.Range("I122").Value = 2790.6365
$ws.Range("K122").Value = 8371.9095
$ws.Range("M122").Value = -5921.9095

$ws.Range("H126").Value = 426772.4
$ws.Range("I126").Value = 719482.7
$ws.Range("J126").Value = 16978
$ws.Range("K126").Value = 2158448.1
$ws.Range("L126").Value = 50934
$ws.Range("M126").Value = -2155978.1
$ws.Range("N126").Value = -55874

$ws.Range("H132").Value = 4890.6943
$ws.Range("I132").Value = 4124.48
$ws.Range("K132").Value = 12373.44
$ws.Range("M132").Value = -9843.439999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 34280.87
$ws.Range("I132").Value = 12494.5
$ws.Range("K132").Value = 37483.5
$ws.Range("M132").Value = -34953.5

$ws.Range("H136").Value = 265361.78
$ws.Range("I136").Value = 251890.92
$ws.Range("K136").Value = 755672.76
$ws.Range("M136").Value = -753122.76
